$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc = "• Develop and simulate ADAS sensors models (camera, radar, lidar, ultrasonic) with in virtual environments.`n• Create and implement test scenarios to validate ADAS features (e.g., AEB, ACC, LKA) Tools like IPG CarMaker.`n• Work closely with vehicle dynamics to ensure accurate integration of dynamics into simulation platforms.`n• Analyze Sensors data and simulation outputs to evaluate features performance and support development.`n• Implement Python programming as required.`n• Participate in Agile development cycles including sprint planning, reviews, and retrospectives`n• Collaborate effectively with team members and stakeholders, maintaining clear and proactive communication`n• Maintain and write clean, scalable, and well-documented code.`n"

$ws.Range("A3").Value = $desc
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 4
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(3).AutoFit()
